# Change the year in the astromap link: 2018 -> 2022.
# The old text is rendered across three differently-formatted runs:
#   "("  +  "http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/" (hyperlink-styled)  +  ")."
# It gets replaced by a single run of plain (unstyled) text containing the
# whole "(url)." string, preceded by a leftover empty run.

$d = $word.ActiveDocument

$oldText = "(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/)."
$newText = "(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

$rng = $d.Content
$found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    $start = $rng.Start
    $end = $rng.End

    # Insert the replacement as a brand-new, unformatted run right before the
    # matched text.
    $d.Range($start, $start).InsertBefore($newText)

    # Leave a (now empty) marker run just ahead of it.
    $d.Range($start, $start).InsertBefore("")

    # Remove the old, differently-formatted runs (the match text got pushed
    # forward by the length of the text we just inserted).
    $d.Range($start + $newText.Length, $end + $newText.Length).Delete()

    Write-Host "Updated astromap link to 2022."
} else {
    Write-Host "Could not find the 2018 astromap link text; no changes made."
}
